$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing
# B:F (values) to C:G, carrying their values/formats along, and
# leaves a new (initially unformatted, copied from column A) column B.
$ws.Columns.Item(2).Insert()

# New header for column B needs the same bold/bordered header style
# that the other header cells (now C1:G1) have, so copy formats over
# from the neighboring header cell.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122) # xlPasteFormats

# The new column B (for rows 2-20) should not carry the bold/border
# style that got copied over from column A during the insert.
$ws.Range("B2:B20").ClearFormats()

# New header for column B
$ws.Range("B1").Value = "segments"

# Segment names, in row order (rows 2-20)
$segments = @(
    "background",
    "back_bumper",
    "back_glass",
    "back_left_door",
    "back_left_light",
    "back_right_door",
    "back_right_light",
    "front_bumper",
    "front_glass",
    "front_left_door",
    "front_left_light",
    "front_right_door",
    "front_right_light",
    "hood",
    "left_mirror",
    "right_mirror",
    "tailgate",
    "trunk",
    "wheel"
)

for ($i = 0; $i -lt $segments.Length; $i++) {
    $row = $i + 2
    # Column A becomes the numeric segment index (0-based)
    $ws.Cells.Item($row, 1).Value = $i
    # Column B gets the segment name text that used to live in column A
    $ws.Cells.Item($row, 2).Value = $segments[$i]
}
